$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Tester / demo1@example / When will the iPad Air be in stock again
$ws.Range("A2").Value = "Tester"
$ws.Range("A2").Font.Color = 0

# Hyperlink on the pre-existing B1 value (demo@email.com), then the newly
# typed B2/B3/B4 email-looking values - mirrors Excel's automatic "convert
# to hyperlink" behaviour as each address is entered. The cell text is set
# first so Hyperlinks.Add keeps it as-is instead of falling back to the
# raw address as the display text.
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:demo@email.com")

$ws.Range("B2").Value = "demo1@example"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:demo1@example")
$ws.Range("C2").Value = "When will the iPad Air be in stock again"

# Row 3: Te / demo1@example.com / When will the iPad Air be in stock again
$ws.Range("A3").Value = "Te"
$ws.Range("B3").Value = "demo1@example.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:demo1@example.com")
$ws.Range("C3").Value = "When will the iPad Air be in stock again"

# Row 4: Tester / demo1@example.com / abcd
$ws.Range("A4").Value = "Tester"
$ws.Range("B4").Value = "demo1@example.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:demo1@example.com")
$ws.Range("C4").Value = "abcd"

# Cosmetic: page set to portrait orientation, active cell moved to E8
$ws.PageSetup.Orientation = 1
$ws.Range("E8").Select()
